$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Some price values (e.g. "228.43") parse as plain numbers, which would
    # make Excel store them as numeric cells instead of text. Forcing the
    # number format to Text ("@") before the assignment keeps the literal
    # string (with leading zeros / exact decimals) intact; switching the
    # style back to Normal afterwards avoids leaving a custom number format
    # applied to the cell (matching the original default-styled text cells).
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "38.226.48"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.092.98"
$ws.Range("E3").Value = "  +2.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "228.43"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.40%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "60.93"
$ws.Range("E7").Value = "  +0.88%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.05%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0851"
$ws.Range("E10").Value = "  +3.88%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.13%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.402.25"
$ws.Range("E12").Value = "  +2.84%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "14.70"
$ws.Range("E13").Value = "  +1.41%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "22.30"
$ws.Range("E14").Value = "  +5.13%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +5.86%  "

# Row 16 - Polygon
Set-TextValue $ws.Range("D16") "0.777"
$ws.Range("E16").Value = "  +2.09%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.104.17"
$ws.Range("E17").Value = "  +3.09%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "38.146.30"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "6.03"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "70.23"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.56%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "224.05"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.10%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.98%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.73%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "169.85"
$ws.Range("E26").Value = "  +1.54%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +1.11%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.12%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "18.98"
$ws.Range("E29").Value = "  +0.52%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +6.46%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -0.45%  "

# Row 32 - WEMIXToken
Set-TextValue $ws.Range("D32") "2.38"
$ws.Range("E32").Value = "  +6.68%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +4.34%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +0.53%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0605"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +4.78%  "

# Row 37 - THORChain
Set-TextValue $ws.Range("D37") "6.41"
$ws.Range("E37").Value = "  +0.91%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +5.46%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  +0.15%  "

# Row 40 - InjectiveProtocol
Set-TextValue $ws.Range("D40") "18.10"
$ws.Range("E40").Value = "  +2.56%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.550.93"
$ws.Range("E41").Value = "  +1.20%  "

# Row 42 - Aave
Set-TextValue $ws.Range("D42") "100.05"
$ws.Range("E42").Value = "  +3.84%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +0.78%  "

# Row 44 - HuobiToken
$ws.Range("E44").Value = "  +1.02%  "

# Row 45 - Cronos
Set-TextValue $ws.Range("D45") "0.0914"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - FTXToken
Set-TextValue $ws.Range("D46") "4.18"
$ws.Range("E46").Value = "  +4.98%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +1.70%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "7.48"
$ws.Range("E48").Value = "  +4.94%  "

# Row 49 - ARBITRUM
Set-TextValue $ws.Range("D49") "1.03"
$ws.Range("E49").Value = "  +1.52%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +0.97%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.288.61"
$ws.Range("E51").Value = "  +2.86%  "
